# Applies the timesheet update described by the commit:
#   - Extends entries for the week of 7/26-8/1 (rows 34-38)
#   - Updates the Week 6 wrap-up note (row 33, C33) and its logged time (B33)
#   - Adds the running-total formula for the new week (E40)
#   - Back-fills the Day column (dates) for the remaining blank weeks (rows 41-56)
#   - Records the two final-due-date notes (rows 55-56)
# Row 59's TOTAL and the weekly SUM formulas in column E recalculate on their own.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 33: finish out "Week 6" row with the longer note + updated hours ---
$ws.Range("B33").Value = 0.125
$ws.Range("C33").Value = "Archives project (Studio Migration) work; begin research for Feedback reading list and create list of sources from HBR, MSLibrary (getAbstract and books)"

# --- Row 34: previously blank except for the date; pick up style from row 35 ---
# (row 34's B cell used a one-off date-style xf that is no longer needed once
#  it carries the same time-of-day value as its neighbours)
$ws.Range("B35").Copy()
$ws.Range("B34").PasteSpecial(-4122)
$ws.Range("B34").Value = 0.14583333333333334
$ws.Range("C34").Value = "Finishing work on User Study white paper & Sync with Kiran; Work on paper/outline for EI Playbook; Feedback reading list"

# --- Rows 35-38: fill in hours + activity notes ---
$ws.Range("B35").Value = 0.125
$ws.Range("C35").Value = "Finish work on User Study white paper; Call with Alex D (Omdia) to discuss agenda/game plan for EI event in Q3 (we discussed possible content areas and topics and he will follow-up with me early next week with Analysts who can participate); Archives project work"

$ws.Range("B36").Value = 0.125
$ws.Range("C36").Value = "Team meeting, work on EI Playbook, Archives project"

$ws.Range("B37").Value = 0.16666666666666666
$ws.Range("C37").Value = "Review of UX Study documentation to prepare for meeting; Team Brainstorming Meeting; work on Feedback Reading list (around 25 sources collected)"

$ws.Range("B38").Value = 0.16666666666666666
$ws.Range("C38").Value = "Sync with Philippe, Expert Insights Playbook work; Draft complete of Feedback Reading List (sent to Kiran for review)"

# --- Row 40: add the week's running-total formula (matches style of E34-E38) ---
$ws.Range("E34").Copy()
$ws.Range("E40").PasteSpecial(-4122)
$ws.Range("E40").Formula = "=SUM(B34:B40)"

# --- Rows 41-47: back-fill the Day column with dates (8/2 - 8/8) ---
$ws.Range("A41").Value = 44410
$ws.Range("A42").Value = 44411
$ws.Range("A43").Value = 44412
$ws.Range("A44").Value = 44413
$ws.Range("A45").Value = 44414
$ws.Range("A46").Value = 44415
$ws.Range("A47").Value = 44416

# --- Rows 48-54: back-fill the Day column with dates (8/9 - 8/15) ---
$ws.Range("A48").Value = 44417
$ws.Range("A49").Value = 44418
$ws.Range("A50").Value = 44419
$ws.Range("A51").Value = 44420
$ws.Range("A52").Value = 44421
$ws.Range("A53").Value = 44422
$ws.Range("A54").Value = 44423

# --- Rows 55-56: final due-date notes ---
$ws.Range("A55").Value = 44424
$ws.Range("C55").Value = "Final Evaluation due for DFW (Philippe)"

$ws.Range("A56").Value = 44429
$ws.Range("C56").Value = "Final Reflection and Artifacts due (Erika)"

# --- Update the view so the active cell / scroll position matches ---
$ws.Range("G59").Select()
